# The author/affiliation metadata for the two OpenAlex records in this sheet
# had been entered on the wrong rows - the "DOM" (row 2) and "Banner" (row 3)
# records need their data swapped so each record's id/title/etc. line up with
# the correct author id formatting.
#
# This swaps the full contents (columns A:Q) of row 2 and row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 17  # Q

$row2Vals = @()
$row3Vals = @()

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $row2Vals += ,$ws.Cells.Item(2, $c).Value()
    $row3Vals += ,$ws.Cells.Item(3, $c).Value()
}

for ($c = $firstCol; $c -le $lastCol; $c++) {
    # Force text formatting so values like dates ("2023-06-01") and
    # numeric-looking strings ("0", "2023") are written back as plain text,
    # matching the original inline-string (text) cell contents instead of
    # being auto-converted to dates/numbers by Excel.
    $ws.Cells.Item(2, $c).NumberFormat = "@"
    $ws.Cells.Item(2, $c).Value = $row3Vals[$c - $firstCol]

    $ws.Cells.Item(3, $c).NumberFormat = "@"
    $ws.Cells.Item(3, $c).Value = $row2Vals[$c - $firstCol]
}
